# Update the "dSF" column (column F) values on the active worksheet.
# This corresponds to a data re-pull for kikuchi_yusei.xlsx where the
# "dSF" figures differ from the previously-pulled "dS0" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = -3
    3  = -3
    5  = -5
    7  = 3
    8  = -7
    11 = 2
    12 = -4
    14 = -6
    15 = -4
    16 = 4
    17 = 4
    18 = 4
    19 = 4
    20 = 3
    21 = 2
    22 = 3
    23 = -3
    24 = 3
    25 = -3
    26 = 1
    27 = -2
    28 = -3
    32 = -1
    33 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 6).Value = $changes[$row]
}
